$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A, shifting existing data
# (currently in A:B) over to B:C.
$ws.Columns("A:A").Insert()

# Narrow the newly inserted column A (closest the host's column-width
# quantization allows to the recorded 3.28515625 OOXML width).
$ws.Columns("A:A").ColumnWidth = 2.5

# Update the selected cell to match the saved view state.
$ws.Range("B11").Select()
